# x1049 addition: support xenium barcode (prebarcode) in file slide reg
#
# Sheet4 gets a new column inserted before the existing "Section address"
# column (old column E) to hold the new "Xenium slide barcode" field.
# The new column inherits the formatting of the column that gets pushed
# to the right (old column E), matches two rows with real values
# ("Xenium slide barcode (...)" header and "ABC1" sample data) and leaves
# the other two data rows blank in the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# Insert a new blank column before column E; everything from E..R shifts to F..S.
$ws.Range("E1").EntireColumn.Insert()

# New column should carry the same formatting that the (now shifted) old
# column E carries in each row - copy format only from column F into E.
$ws.Range("F2:F5").Copy()
$ws.Range("E2:E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new column's content.
$ws.Range("E2").Value = "Xenium slide barcode (…)"
$ws.Range("E4").Value = "ABC1"

# Leave E3 and E5 empty (they stay blank, matching the rest of that row's data gap).

$ws.Range("E5").Select() | Out-Null
